$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.039.95"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "3.511.03"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("E9").Value = "  +6.59%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "4.118.14"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "67.028.88"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "3.503.33"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("E25").Value = "  -4.57%  "
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.895"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0745"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "2.794.52"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0306"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.848"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
